$wb = $excel.ActiveWorkbook

# ---- Sheet: Matriz_Resultados ----
$ws = $wb.Worksheets.Item("Matriz_Resultados")
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("B9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0

# ---- Sheet: P_valores ----
$ws = $wb.Worksheets.Item("P_valores")
$ws.Range("C2").Value = 0.009348655955216501
$ws.Range("D2").Value = 0.0145739921536836
$ws.Range("E2").Value = 0.03895749117734137
$ws.Range("F2").Value = 0.02097010059017079
$ws.Range("G2").Value = 0.008770594071426219
$ws.Range("H2").Value = 0.01041051977450214
$ws.Range("I2").Value = 0.008602401169364837
$ws.Range("J2").Value = 0.1805340155793134
$ws.Range("B3").Value = 0.009348655955216501
$ws.Range("D3").Value = 0.00004409311725406262
$ws.Range("E3").Value = 0.001285292206323785
$ws.Range("F3").Value = 0.009264341909712037
$ws.Range("G3").Value = 0.01763273101747753
$ws.Range("H3").Value = 0.01129677035920218
$ws.Range("I3").Value = 0.03814544533403308
$ws.Range("J3").Value = 0.0004542018997495401
$ws.Range("B4").Value = 0.0145739921536836
$ws.Range("C4").Value = 0.00004409311725406262
$ws.Range("E4").Value = 0.003402129693828693
$ws.Range("F4").Value = 0.01524936719278558
$ws.Range("G4").Value = 0.09264238494183585
$ws.Range("H4").Value = 0.07836846013718857
$ws.Range("I4").Value = 0.6414060685503875
$ws.Range("J4").Value = 0.0008141827819834191
$ws.Range("B5").Value = 0.03895749117734137
$ws.Range("C5").Value = 0.001285292206323785
$ws.Range("D5").Value = 0.003402129693828693
$ws.Range("F5").Value = 0.05045272604095907
$ws.Range("G5").Value = 0.02361361202704138
$ws.Range("H5").Value = 0.006130563304465486
$ws.Range("I5").Value = 0.005922490574846284
$ws.Range("J5").Value = 0.001359964783045342
$ws.Range("B6").Value = 0.02097010059017079
$ws.Range("C6").Value = 0.009264341909712037
$ws.Range("D6").Value = 0.01524936719278558
$ws.Range("E6").Value = 0.05045272604095907
$ws.Range("G6").Value = 0.008221848978701152
$ws.Range("H6").Value = 0.01041599587972541
$ws.Range("I6").Value = 0.008399239234146316
$ws.Range("J6").Value = 0.3282394305669851
$ws.Range("B7").Value = 0.008770594071426219
$ws.Range("C7").Value = 0.01763273101747753
$ws.Range("D7").Value = 0.09264238494183585
$ws.Range("E7").Value = 0.02361361202704138
$ws.Range("F7").Value = 0.008221848978701152
$ws.Range("H7").Value = 0.3993389294435115
$ws.Range("I7").Value = 0.07494986794828629
$ws.Range("J7").Value = 0.0002882213082229246
$ws.Range("B8").Value = 0.01041051977450214
$ws.Range("C8").Value = 0.01129677035920218
$ws.Range("D8").Value = 0.07836846013718857
$ws.Range("E8").Value = 0.006130563304465486
$ws.Range("F8").Value = 0.01041599587972541
$ws.Range("G8").Value = 0.3993389294435115
$ws.Range("I8").Value = 0.09409883306369204
$ws.Range("J8").Value = 0.0002072739114911126
$ws.Range("B9").Value = 0.008602401169364837
$ws.Range("C9").Value = 0.03814544533403308
$ws.Range("D9").Value = 0.6414060685503875
$ws.Range("E9").Value = 0.005922490574846284
$ws.Range("F9").Value = 0.008399239234146316
$ws.Range("G9").Value = 0.07494986794828629
$ws.Range("H9").Value = 0.09409883306369204
$ws.Range("J9").Value = 0.0002667018596558268
$ws.Range("B10").Value = 0.1805340155793134
$ws.Range("C10").Value = 0.0004542018997495401
$ws.Range("D10").Value = 0.0008141827819834191
$ws.Range("E10").Value = 0.001359964783045342
$ws.Range("F10").Value = 0.3282394305669851
$ws.Range("G10").Value = 0.0002882213082229246
$ws.Range("H10").Value = 0.0002072739114911126
$ws.Range("I10").Value = 0.0002667018596558268

# ---- Sheet: Estadisticos_DM ----
$ws = $wb.Worksheets.Item("Estadisticos_DM")
$ws.Range("C2").Value = 2.848307454412941
$ws.Range("D2").Value = 2.651609509664682
$ws.Range("E2").Value = 2.195640238823384
$ws.Range("F2").Value = 2.486627551202071
$ws.Range("G2").Value = 2.876226202826496
$ws.Range("H2").Value = 2.80105814201709
$ws.Range("I2").Value = 2.884679608026738
$ws.Range("J2").Value = 1.383025466157474
$ws.Range("B3").Value = -2.848307454412941
$ws.Range("D3").Value = -5.07358700442393
$ws.Range("E3").Value = -3.688513177280016
$ws.Range("F3").Value = -2.852275286457594
$ws.Range("G3").Value = -2.565699054335158
$ws.Range("H3").Value = -2.765008918929533
$ws.Range("I3").Value = -2.205786134408974
$ws.Range("J3").Value = -4.116260530357766
$ws.Range("B4").Value = -2.651609509664682
$ws.Range("C4").Value = 5.07358700442393
$ws.Range("E4").Value = -3.282267999678816
$ws.Range("F4").Value = -2.63127198169974
$ws.Range("G4").Value = -1.758067872748599
$ws.Range("H4").Value = -1.846157739044234
$ws.Range("I4").Value = -0.4722401901876879
$ws.Range("J4").Value = -3.876766447159473
$ws.Range("B5").Value = -2.195640238823384
$ws.Range("C5").Value = 3.688513177280016
$ws.Range("D5").Value = 3.282267999678816
$ws.Range("F5").Value = -2.069417690080629
$ws.Range("G5").Value = 2.431926170549513
$ws.Range("H5").Value = 3.0314406922746
$ws.Range("I5").Value = 3.046288351131652
$ws.Range("J5").Value = -3.665149033151636
$ws.Range("B6").Value = -2.486627551202071
$ws.Range("C6").Value = 2.852275286457594
$ws.Range("D6").Value = 2.63127198169974
$ws.Range("E6").Value = 2.069417690080629
$ws.Range("G6").Value = 2.904404162684012
$ws.Range("H6").Value = 2.800826576384943
$ws.Range("I6").Value = 2.895103589286281
$ws.Range("J6").Value = 0.9998812944066217
$ws.Range("B7").Value = -2.876226202826496
$ws.Range("C7").Value = 2.565699054335158
$ws.Range("D7").Value = 1.758067872748599
$ws.Range("E7").Value = -2.431926170549513
$ws.Range("F7").Value = -2.904404162684012
$ws.Range("H7").Value = 0.8594905532491552
$ws.Range("I7").Value = 1.869310308183943
$ws.Range("J7").Value = -4.302446226804591
$ws.Range("B8").Value = -2.80105814201709
$ws.Range("C8").Value = 2.765008918929533
$ws.Range("D8").Value = 1.846157739044234
$ws.Range("E8").Value = -3.0314406922746
$ws.Range("F8").Value = -2.800826576384943
$ws.Range("G8").Value = -0.8594905532491552
$ws.Range("I8").Value = 1.7497514818427
$ws.Range("J8").Value = -4.437369341256422
$ws.Range("B9").Value = -2.884679608026738
$ws.Range("C9").Value = 2.205786134408974
$ws.Range("D9").Value = 0.4722401901876879
$ws.Range("E9").Value = -3.046288351131652
$ws.Range("F9").Value = -2.895103589286281
$ws.Range("G9").Value = -1.869310308183943
$ws.Range("H9").Value = -1.7497514818427
$ws.Range("J9").Value = -4.334200029440385
$ws.Range("B10").Value = -1.383025466157474
$ws.Range("C10").Value = 4.116260530357766
$ws.Range("D10").Value = 3.876766447159473
$ws.Range("E10").Value = 3.665149033151636
$ws.Range("F10").Value = -0.9998812944066217
$ws.Range("G10").Value = 4.302446226804591
$ws.Range("H10").Value = 4.437369341256422
$ws.Range("I10").Value = 4.334200029440385

# ---- Sheet: Resumen ----
$ws = $wb.Worksheets.Item("Resumen")
$ws.Range("B2").Value = 3
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 37.5
$ws.Range("B3").Value = 1
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 12.5
$ws.Range("B4").Value = 1
$ws.Range("D4").Value = 6
$ws.Range("E4").Value = 12.5
$ws.Range("A5").Value = "LSPMW"
$ws.Range("B5").Value = 1
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 12.5
$ws.Range("F5").Value = 1.40982501785926
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 12.5
$ws.Range("A7").Value = "AV-MCPS"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 7
$ws.Range("E7").Value = 12.5
$ws.Range("F7").Value = 1.070183749104134
$ws.Range("A8").Value = "Block Bootstrapping"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 2.468373163546003
$ws.Range("A9").Value = "AREPD"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 8
$ws.Range("F9").Value = 2.270482596061766
